$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 231, shifting rows 231:355 down to 232:356
$ws.Rows.Item(231).Insert()

# Fill the new row 231 with its data (same as original row 231 except a few changed fields)
$ws.Range("A231").Value = 6
$ws.Range("B231").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C231").Value = "Metropolitana"
$ws.Range("D231").Value = 44572
$ws.Range("D231").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E231").Value = 13
$ws.Range("F231").Value = 100112039
$ws.Range("G231").Value = "Ciboulette"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 1070
$ws.Range("K231").Value = 800
$ws.Range("L231").Value = 900
$ws.Range("M231").Value = 841
$ws.Range("N231").Value = "`$/docena de atados"
$ws.Range("O231").Value = "Provincia de Quillota"
$ws.Range("P231").Value = 280
$ws.Range("Q231").Value = 3
$ws.Range("R231").Value = "Hortaliza"
